$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B-column accuracy values (per-epoch accuracy changed on re-run) ---
$ws.Range("B4").Value = 0.953125
$ws.Range("B5").Value = 0.921875
$ws.Range("B6").Value = 0.90625
$ws.Range("B7").Value = 0.9375
$ws.Range("B8").Value = 0.90625
$ws.Range("B9").Value = 0.875
$ws.Range("B11").Value = 0.921875
$ws.Range("B13").Value = 0.890625
$ws.Range("B14").Value = 0.859375
$ws.Range("B15").Value = 0.859375
$ws.Range("B16").Value = 0.875
$ws.Range("B17").Value = 0.890625
$ws.Range("B18").Value = 0.859375
$ws.Range("B19").Value = 0.875
$ws.Range("B20").Value = 0.859375
$ws.Range("B21").Value = 0.875
$ws.Range("B22").Value = 0.890625
$ws.Range("B23").Value = 0.875
$ws.Range("B24").Value = 0.875
$ws.Range("B25").Value = 0.875
$ws.Range("B26").Value = 0.890625
$ws.Range("B27").Value = 0.875
$ws.Range("B28").Value = 0.90625
$ws.Range("B29").Value = 0.875
$ws.Range("B30").Value = 0.859375
$ws.Range("B31").Value = 0.84375
$ws.Range("B32").Value = 0.8125
$ws.Range("B33").Value = 0.828125
$ws.Range("B34").Value = 0.875
$ws.Range("B35").Value = 0.859375
$ws.Range("B37").Value = 0.8125
$ws.Range("B38").Value = 0.84375
$ws.Range("B39").Value = 0.859375
$ws.Range("B40").Value = 0.859375
$ws.Range("B41").Value = 0.875
$ws.Range("B42").Value = 0.859375
$ws.Range("B43").Value = 0.859375
$ws.Range("B44").Value = 0.84375
$ws.Range("B45").Value = 0.8125
$ws.Range("B46").Value = 0.875
$ws.Range("B47").Value = 0.875
$ws.Range("B48").Value = 0.859375
$ws.Range("B49").Value = 0.84375
$ws.Range("B50").Value = 0.84375
$ws.Range("B51").Value = 0.828125
$ws.Range("B52").Value = 0.828125
$ws.Range("B53").Value = 0.828125
$ws.Range("B54").Value = 0.859375
$ws.Range("B57").Value = 0.84375
$ws.Range("B58").Value = 0.859375
$ws.Range("B60").Value = 0.828125
$ws.Range("B61").Value = 0.828125
$ws.Range("B62").Value = 0.84375
$ws.Range("B63").Value = 0.8125
$ws.Range("B64").Value = 0.875
$ws.Range("B65").Value = 0.859375
$ws.Range("B66").Value = 0.875
$ws.Range("B67").Value = 0.875
$ws.Range("B68").Value = 0.875
$ws.Range("B69").Value = 0.875
$ws.Range("B70").Value = 0.859375
$ws.Range("B71").Value = 0.875
$ws.Range("B72").Value = 0.875
$ws.Range("B73").Value = 0.890625
$ws.Range("B74").Value = 0.890625
$ws.Range("B75").Value = 0.859375
$ws.Range("B76").Value = 0.859375
$ws.Range("B77").Value = 0.859375
$ws.Range("B78").Value = 0.859375
$ws.Range("B79").Value = 0.875
$ws.Range("B80").Value = 0.859375
$ws.Range("B81").Value = 0.875
$ws.Range("B82").Value = 0.84375
$ws.Range("B83").Value = 0.859375
$ws.Range("B84").Value = 0.859375
$ws.Range("B85").Value = 0.859375
$ws.Range("B86").Value = 0.859375
$ws.Range("B87").Value = 0.859375
$ws.Range("B88").Value = 0.859375
$ws.Range("B89").Value = 0.875
$ws.Range("B90").Value = 0.875
$ws.Range("B91").Value = 0.890625
$ws.Range("B92").Value = 0.875
$ws.Range("B93").Value = 0.875
$ws.Range("B94").Value = 0.875
$ws.Range("B95").Value = 0.875
$ws.Range("B96").Value = 0.875
$ws.Range("B97").Value = 0.875
$ws.Range("B98").Value = 0.875
$ws.Range("B99").Value = 0.875
$ws.Range("B100").Value = 0.875
$ws.Range("B101").Value = 0.875
$ws.Range("B102").Value = 0.875
$ws.Range("B103").Value = 0.78125
$ws.Range("B104").Value = 0.84375
$ws.Range("B105").Value = 0.828125
$ws.Range("B106").Value = 0.796875
$ws.Range("B107").Value = 0.78125
$ws.Range("B109").Value = 0.7058823529411765

# --- Update A102:A109 repr strings (object memory address changed on re-run) ---
$ws.Range("A102:A109").Value = "<__main__.DisplayOutputs object at 0x7fcf9fb1fe50>"

# --- Update the saved sheet selection/view state ---
$ws.Range("A1:XFD1048576").Select()
$ws.Range("S21").Activate()
